# CPT_QTR_FIN.xlsx quarterly refresh:
# Two new quarter columns (most-recent first) are inserted immediately
# before the existing "D" column, pushing the old D:K data right to F:M.
# The two freshly inserted columns (new D:E) are then populated with the
# latest two quarters of data for every reported line item.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert two new blank columns at D:E -- shifts existing D:K -> F:M.
$ws.Range("D:E").EntireColumn.Insert()

# 2) Style for the two new date-header rows (Period Ending).
$dateFormat = "[$-409]d\-mmm\-yy;@"
# Style used by every ordinary data row (matches the pre-existing "#,##0" xf).
$numFormat = "#,##0"

$dateRows = @(7, 38, 80)

# 3) New values for the two inserted columns, keyed by row number.
#    Each entry is the (D, E) pair -- D is the newest quarter, E the one
#    before it. $null means "leave the cell blank" (spacer rows).
$newValues = @{
    7   = @(43465, 43373)
    8   = @(244900, 241800)
    9   = @(94100, 95300)
    10  = @(150800, 146500)
    11  = @($null, $null)
    12  = @("NA", "NA")
    13  = @(0, 0)
    14  = @(0, 0)
    15  = @(78700, 76500)
    16  = @($null, $null)
    17  = @(206500, 203400)
    18  = @(38400, 38400)
    19  = @($null, $null)
    20  = @(2200, 1900)
    21  = @(119300, 116800)
    22  = @(0, 0)
    23  = @(40600, 40300)
    24  = @(300, 300)
    25  = @(0, 0)
    26  = @(40300, 40000)
    27  = @(38900, 38600)
    28  = @(0, 0)
    29  = @(0, "NA")
    30  = @(0, 0)
    31  = @(0, 0)
    32  = @(-2200, -1900)
    33  = @(38900, 38600)
    34  = @(0, 0)
    35  = @(38900, 38600)
    38  = @(43465, 43373)
    39  = @($null, $null)
    40  = @($null, $null)
    41  = @(34400, 8500)
    42  = @(0, 0)
    43  = @(22900, 22600)
    44  = @(0, 0)
    45  = @(0, 0)
    46  = @(0, 0)
    47  = @(22300, 24700)
    48  = @(5925300, 5904200)
    49  = @(0, 0)
    50  = @(0, 0)
    51  = @(0, 0)
    52  = @(9200, 10100)
    53  = @(0, 0)
    54  = @(6219600, 6198500)
    55  = @($null, $null)
    56  = @($null, $null)
    57  = @(146900, 140000)
    58  = @(0, 0)
    59  = @(129300, 145200)
    60  = @(0, 0)
    61  = @(2321600, 2259600)
    62  = @(0, 0)
    63  = @(0, 0)
    64  = @(0, 0)
    65  = @(0, 0)
    66  = @(2908200, 2858500)
    67  = @($null, $null)
    68  = @(0, 0)
    69  = @(0, 0)
    70  = @(0, 0)
    71  = @(0, 0)
    72  = @(-495500, -466500)
    73  = @(0, 0)
    74  = @(0, 0)
    75  = @(0, 0)
    76  = @(3311400, 3340000)
    77  = @(0, 0)
    80  = @(43465, 43373)
    81  = @(38900, 38600)
    82  = @($null, $null)
    83  = @(78700, 76500)
    84  = @(0, 0)
    85  = @(0, 0)
    86  = @(0, 0)
    87  = @(0, 0)
    88  = @(0, 0)
    89  = @(127800, 144700)
    90  = @($null, $null)
    91  = @(-86900, -94800)
    92  = @(0, 0)
    93  = @(0, 0)
    94  = @(-89100, -164600)
    95  = @($null, $null)
    96  = @(-75000, -75100)
    97  = @(0, 0)
    98  = @(0, 0)
    99  = @(0, 0)
    100 = @(-13700, -35200)
    101 = @(0, 0)
    102 = @(25000, -55100)
}

foreach ($r in $newValues.Keys) {
    $pair = $newValues[$r]
    $rng = $ws.Range("D" + $r + ":E" + $r)

    if ($dateRows -contains $r) {
        $rng.NumberFormat = $dateFormat
    } else {
        $rng.NumberFormat = $numFormat
    }

    $dCell = $ws.Cells.Item($r, 4)
    $eCell = $ws.Cells.Item($r, 5)

    if ($pair[0] -ne $null) { $dCell.Value = $pair[0] }
    if ($pair[1] -ne $null) { $eCell.Value = $pair[1] }
}
